$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5 previously held the number 5; it now holds the text "testxxxxx"
$ws.Range("D5").Value = "testxxxxx"

# A19 / A20 previously held the number 1; they now hold the text "testxxxxx"
$ws.Range("A19").Value = "testxxxxx"
$ws.Range("A20").Value = "testxxxxx"

# New row 21 with the same text value in column A
$ws.Range("A21").Value = "testxxxxx"

# Move / record the active selection on the newly added cell
[void]$ws.Range("A21").Select()
